$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Locate the "Why are ethics important with AI?" heading paragraph
#    and the paragraph right after it that only holds the _GoBack
#    bookmark (empty paragraph with bookmarkStart/bookmarkEnd).
# ------------------------------------------------------------------
$headingIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $txt = $p.Range.Text
    if ($txt -eq "Why are ethics important with AI?`r") {
        $headingIndex = $i
    }
}

$heading = $d.Paragraphs($headingIndex)

# ------------------------------------------------------------------
# 2) Move the _GoBack bookmark so it sits at the end of the heading
#    paragraph (after the existing run, still inside that <w:p>)
#    instead of living in its own paragraph.
# ------------------------------------------------------------------
$tail = $heading.Range.Duplicate
[void]$tail.MoveEnd(1, -1)
$tail.Collapse(0)
$tail.InsertAfter("~")
$d.Bookmarks.Add("_GoBack", $tail)

# remove the placeholder character again, the bookmark stays behind
$placeholder = $d.Range($tail.Start, $tail.End)
$placeholder.Delete()

# ------------------------------------------------------------------
# 3) Turn the (now bookmark-free) empty paragraph into the new
#    "In a world where..." paragraph.
# ------------------------------------------------------------------
$newPara = $d.Paragraphs($headingIndex + 1)
$newPara.Range.Text = "In a world where humans are becoming more and more dependent on machines, the need for AI is exponentially increasing. Due to this circumstance, ensuring that any machines or computer used by machines that are safety critical is "

# ------------------------------------------------------------------
# 4) Drop the stray <w:lastRenderedPageBreak/> in front of "Notes".
# ------------------------------------------------------------------
$j = 0
foreach ($p in $d.Paragraphs) {
    $j = $j + 1
    if ($p.Range.Text -eq "Notes`r") {
        $notesRange = $p.Range.Duplicate
        [void]$notesRange.MoveEnd(1, -1)
        $notesRange.Delete()
        $notesRange.InsertAfter("Notes")
    }
}
